$d = $word.ActiveDocument

# Ativação date update
$d.Content.Find.Execute('Ativação: 01/01/2018', $true, $false, $false, $false, $false, $true, 1, $false, 'Ativação: 01/01/2025', 2) | Out-Null

# Objetivos (PT) - full paragraph replace
$d.Content.Find.Execute('Desenvolver o aprendizado teórico e prático da Bioquímica através da execução de práticas de laboratório baseadas na evolução do conteúdo teórico ministrado na disciplina Bioquímica II.', $true, $false, $false, $false, $false, $true, 1, $false, 'Desenvolver e aperfeiçoar o entendimento teórico e prático dos processos bioquímicos fundamentais através da realização de atividades práticas de laboratório.', 2) | Out-Null

# Objetivos (EN) - insert text into the empty italic run
$pObjEn = $d.Paragraphs.Item(7)
$pObjEn.Range.InsertBefore('Developing and enhancing the theoretical and practical understanding of fundamental biochemical processes through the execution of laboratory practical activities.')

# Programa resumido (PT) - full paragraph replace
$d.Content.Find.Execute('Propriedades gerais de glicídios; Fermentação anaeróbia; Extração deClorofila e Reação de Hill; Transporte de glicídios e indução de enzimas.', $true, $false, $false, $false, $false, $true, 1, $false, 'Reação de saponificação; Enzimas proteolíticas em produtos comerciais; Extração líquido-líquido de proteínas; Biomateriais sustentáveis; Produção e destilação de etanol; Precipitação de biomoléculas e Reação de Hill.', 2) | Out-Null

# Programa resumido (EN) - paragraph 12, direct range text (source text contains a non-standard character)
$pResumoEn = $d.Paragraphs.Item(12)
$pResumoEn.Range.Text = 'Saponification reaction; Proteolytic enzymes in commercial products; Liquid-liquid extraction of proteins; Sustainable biomaterials; Production and distillation of ethanol; Precipitation of biomolecules; and Hill reaction.'

# Programa (PT) - paragraph 14, direct range text (new text contains straight quotes)
$pProgPt = $d.Paragraphs.Item(14)
$pProgPt.Range.Text = 'Aplicação da saponificação em processos industriais, agentes envolvidos na reação e sua utilização em produtos comerciais (cálculo de rendimento, CMC e pH). Avaliação enzimática de proteases para determinação de sua atividade proteolítica em produtos comerciais (sabão em pó, detergentes e cosméticos). Extração líquido-líquido de proteínas e enzimas utilizando solventes orgânicos e polímeros/tensoativos - quantificação dos parâmetros de extração (balanço de massa, recuperação, fator de purificação). Obtenção de biomateriais (bioplástico) de interesse biotecnológico derivado de fontes biológicas - cálculo do rendimento; caracterização do produto final obtido (textura, cor e cheiro) e comparação com os plásticos convencionais. Produção e destilação de etanol - conceitos gerais e fermentação de glicose; produção de etanol e CO2; consumo da fonte de carbono; cálculo da eficiência do processo; ação de um inibidor da glicólise. Precipitação de biomoléculas utilizando diferentes agentes precipitadores (sais, polímeros e solventes orgânicos) - quantificação da recuperação, pH e potencial elétrico. Extração de clorofila e reação de Hill - estrutura de cloroplastos; papel da clorofila nos sistemas fotossintéticos I e II; fase escura/luminosa; produção de NADP; produção de ATP; papel do corante como aceptor de prótons e elétrons. *Dentro do programa da disciplina é planejado realizar eventual "Visita Didática Complementar".'

# Programa (EN) - full paragraph replace
$d.Content.Find.Execute('General proprieties of glycides: main qualitative tests for identification and differentiation of glycides; important reactions and spectrophotometric dosage of reducing monosaccharides.Anaerobic fermentation: general concepts and glucoseFermentation; ethanol and CO2 production; consumption of carbon source; calculation of the process efficiency; action of inhibitor of glycolysis.Chlorophyll extraction andHill reaction: structure of chloroplasts, chlorophyll role in the photosystems I and II; dark phase and light phase; NADP production; dye role as an acceptor ofprotons and electrons. Glycides transportation and enzyme induction: general concepts; enzymes of galactose catabolism; catabolic repression, inactivation and modification; constitutive and induced enzymatic systems in yeast cells.', $true, $false, $false, $false, $false, $true, 1, $false, 'Application of saponification in industrial processes, agents involved in the reaction, and its use in commercial products (yield calculation, CMC, and pH). Enzymatic evaluation of proteases to determine their proteolytic activity in commercial products (laundry detergent, detergents, and cosmetics). Liquid-liquid extraction of proteins and enzymes using organic solvents and polymers/surfactants - quantification of extraction parameters (mass balance, recovery, purification factor). Production of biomaterials (bioplastic) of biotechnological interest derived from biological sources - yield calculation; characterization of the final product obtained (texture, color, and odor) and comparison with conventional plastics. Production and distillation of ethanol - general concepts and glucose fermentation; ethanol and CO2 production; carbon source consumption; process efficiency calculation; action of a glycolysis inhibitor. Precipitation of biomolecules using different precipitating agents (salts, polymers, and organic solvents) - recovery quantification, pH, and electrical potential. Chlorophyll extraction and Hill reaction - chloroplast structure; role of chlorophyll in photosynthetic systems I and II; dark/light phase; NADP production; ATP production; dye role as proton and electron acceptor. *Complementary didactic visit is planned within the course program.', 2) | Out-Null

# Avaliação - Método text run
$d.Content.Find.Execute('A avaliação será feita por meio de uma prova escrita e notas de relatórios (R).', $true, $false, $false, $false, $false, $true, 1, $false, 'A avaliação será realizada através de uma prova escrita (P) e um relatório de atividades práticas (R).', 2) | Out-Null

# Avaliação - Critério text run
$d.Content.Find.Execute('A Nota final (NF) será calculada da seguinte maneira: NF = (P1*2 + R)/3.', $true, $false, $false, $false, $false, $true, 1, $false, 'A nota final (NF) será calculada conforme: NF = (P + R)/2. A', 2) | Out-Null

# Avaliação - Norma de recuperação text run
$d.Content.Find.Execute('A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada pela fórmula: MR = (NF + PR)/2', $true, $false, $false, $false, $false, $true, 1, $false, 'A recuperação será realizada através de uma prova escrita (PR) e a média de recuperação (MR) será calculada conforme: MR = (NF + PR)/2.', 2) | Out-Null

# Bibliografia - full paragraph replace
$d.Content.Find.Execute('CISTERNAS, J. R. Fundamentos de bioquímica experimental. São Paulo : Atheneu, 2005. ISBN: 9788573791075.NELSON, D. L., COX. M. M. Princípios de bioquímica de Lehninger. Porto Alegre : Artmed, 2011. ISBN: 9788536324180.VOET, D., VOET, J. G. Bioquímica. Porto Alegre : Artmed, 2013. ISBN: 9788582710043.', $true, $false, $false, $false, $false, $true, 1, $false, '1. Rocha Filho, J.A., Vitolo, M. Guia para aulas práticas de biotecnologia de enzimas e fermentação. Editora Blucher, 2021. 2. Cisternas, J.R. Fundamentos de bioquímica experimental. São Paulo: Atheneu, 2005. 3. Nelson, D.L., Cox, M.M. Princípios de bioquímica de Lehninger. Artmed Editora, 2022. 4. Voet, D., Voet, J.G., Pratt, C.W. Fundamentos de Bioquímica: a vida em nivel molecular. Artmed Editora, 2014. 5. Vitolo, M., Pessoa Junior, A., Monteiro, G., Carvalho, J.C.M., Stephano, M.A., Sato, S. Biotecnologia farmacêutica: aspectos sobre aplicação industrial. Editora Blucher, 2015.', 2) | Out-Null

Write-Host "Done applying edits"